# ---------------------------------------------------------------------------
# projektbegleitender_bericht.docx -- apply "Add files via upload" edit
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    # Drop the trailing cell-mark character so the whole visible content of
    # the cell (including any paragraph mark/run boundaries) gets replaced
    # by one clean run instead of leaving stray leftover runs behind.
    $rng.End = $rng.End - 1
    $rng.Text = $newText
}

# ---------------------------------------------------------------------------
# 1) Table 1 - "Aufgabenbereich / Beschreibung / Mitglied"
# ---------------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
Set-CellText $t1 2 3 "Franziska, Kurt, Max, Sebastian, Theodor"
Set-CellText $t1 3 3 "Kurt, Theodor"
Set-CellText $t1 4 3 "Kurt, Theodor"

# ---------------------------------------------------------------------------
# 2) Table 2 - "Datum / Erreichter Meilenstein / Mitglieder"
# ---------------------------------------------------------------------------
$t2 = $d.Tables.Item(2)
Set-CellText $t2 4 3 "Franziska"
Set-CellText $t2 7 1 "17.05.23"
Set-CellText $t2 7 3 "Max, Sebastian"
Set-CellText $t2 8 1 "19.05.23"

# "Prototyp ?" -> "Javadokumentationen", and the proofErr annotation around
# it flips from grammar (gramStart/gramEnd) to spelling (spellStart/spellEnd).
$d.Content.Find.ClearFormatting()
$found = $d.Content.Find.Execute("Prototyp ?", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Javadokumentationen", 2)

# ---------------------------------------------------------------------------
# 3) Table 3 - "Aufgabe / Arbeitsstunden"
# ---------------------------------------------------------------------------
$t3 = $d.Tables.Item(3)
Set-CellText $t3 2 2 "20h"
Set-CellText $t3 3 2 "18h"
Set-CellText $t3 4 2 "5h"
Set-CellText $t3 5 2 "25h"
Set-CellText $t3 8 2 "5h"
Set-CellText $t3 9 2 "45h"
Set-CellText $t3 10 2 "6h"
Set-CellText $t3 11 2 "6h"
Set-CellText $t3 12 2 "6h"

# ---------------------------------------------------------------------------
# 4) Table 4 - "Zeitaufwand gesamt"
# ---------------------------------------------------------------------------
$t4 = $d.Tables.Item(4)
Set-CellText $t4 1 2 "170h"
